# Generate Report for Archive
#
# The files "43eebd46-d413-44ca-af25-a940c37a081d.md" (row 3) and
# "891b2b37-e62f-441f-a054-c491cb0d72e9.md" (row 4) have moved out of the
# "Ready for handoff" state into "In Translation" on both the zh-cn and
# de-de localization status sheets. Update the Status column (column C)
# for those two rows accordingly.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"
